$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Step 1: materialize the full six header/footer slots -----------------
# Touching any header/footer slot other than the existing "primary" one
# causes the host to mint header2/header3/footer1/footer2/footer3 and wire
# up the section's six w:headerReference/w:footerReference entries
# (even/default/first x header/footer). The pre-existing default header
# (current header1.xml) is copied forward to become the new "default"
# slot (Headers(1) -> header2.xml); rId1/rId2 hyperlink relationships
# travel with it.
$sec.Footers(3).Range.Text = ""

# --- Step 2: empty out the even-page and first-page headers/footers -------
# wdHeaderFooterEvenPages = 3, wdHeaderFooterFirstPage = 2,
# wdHeaderFooterPrimary   = 1
$sec.Headers(3).Range.Text = ""   # even header  -> header1.xml
$sec.Headers(2).Range.Text = ""   # first header -> header3.xml
$sec.Footers(1).Range.Text = ""   # default footer -> footer2.xml
$sec.Footers(2).Range.Text = ""   # first footer   -> footer3.xml
$sec.Footers(3).Range.Text = ""   # even footer    -> footer1.xml

# --- Step 3: rewrite the default header's "Portfolio" hyperlink as a ------
# HYPERLINK field (begin/instrText/separate/end) instead of a real
# w:hyperlink relationship run.
$hdr = $sec.Headers(1)
$full = $hdr.Range
$xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="0521B427" w14:textId="57DA7515" w:rsidR="00674771" w:rsidRDefault="00674771"><w:pPr><w:pStyle w:val="Header"/></w:pPr><w:r w:rsidRPr="00674771"><w:rPr><w:rFonts w:ascii="Calibri Light" w:eastAsia="Times New Roman" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:kern w:val="0"/><w:sz w:val="18"/><w:szCs w:val="18"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Nathaniel Bryce Carroll</w:t></w:r><w:r w:rsidRPr="00674771"><w:rPr><w:rFonts w:ascii="Calibri Light" w:eastAsia="Times New Roman" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:kern w:val="0"/><w:sz w:val="18"/><w:szCs w:val="18"/><w14:ligatures w14:val="none"/></w:rPr><w:br/></w:r><w:hyperlink r:id="rId1" w:history="1"><w:r w:rsidRPr="00674771"><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Calibri Light" w:eastAsia="Times New Roman" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:kern w:val="0"/><w:sz w:val="18"/><w:szCs w:val="18"/><w14:ligatures w14:val="none"/></w:rPr><w:t>nathaniel.b.carroll@outlook.com</w:t></w:r></w:hyperlink><w:r w:rsidRPr="00674771"><w:rPr><w:rFonts w:ascii="Calibri Light" w:eastAsia="Times New Roman" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:kern w:val="0"/><w:sz w:val="18"/><w:szCs w:val="18"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve"> &#8729; (512) 656-1997 &#8729; [</w:t></w:r><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText>HYPERLINK "https://nathaniel-b-carroll.github.io/" \l "home"</w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Calibri Light" w:eastAsia="Times New Roman" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:kern w:val="0"/><w:sz w:val="18"/><w:szCs w:val="18"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Portfolio</w:t></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Calibri Light" w:eastAsia="Times New Roman" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:kern w:val="0"/><w:sz w:val="18"/><w:szCs w:val="18"/><w14:ligatures w14:val="none"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri Light" w:eastAsia="Times New Roman" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:kern w:val="0"/><w:sz w:val="18"/><w:szCs w:val="18"/><w14:ligatures w14:val="none"/></w:rPr><w:t>] &#8729; Austin, TX</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$full.InsertXML($xml)
